$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 100
$ws.Range("B100").Value = 6732836
$ws.Range("E100").Value = "FK Siauliai"
$ws.Range("F100").Value = "Banga Gargzdai"
$ws.Range("G100").Value = 3
$ws.Range("H100").Value = 0
$ws.Range("I100").Value = 2
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = "H"
$ws.Range("L100").Value = 1.222
$ws.Range("M100").Value = 5.5
$ws.Range("N100").Value = 9
$ws.Range("O100").Value = 1.363
$ws.Range("P100").Value = 4.5
$ws.Range("Q100").Value = 7
$ws.Range("R100").Value = -1.25
$ws.Range("S100").Value = 1.9
$ws.Range("T100").Value = 1.9
$ws.Range("U100").Value = 2.5
$ws.Range("V100").Value = 1.975
$ws.Range("W100").Value = 1.825
$ws.Range("X100").Value = 0.363
$ws.Range("Y100").Value = -1
$ws.Range("Z100").Value = -1
$ws.Range("AA100").Value = 0.8999999999999999
$ws.Range("AB100").Value = -1
$ws.Range("AC100").Value = 0.9750000000000001
$ws.Range("AD100").Value = -1

# Row 101
$ws.Range("B101").Value = 6732834
$ws.Range("E101").Value = "Panevezys"
$ws.Range("F101").Value = "FK Dziugas Telsiai"
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 0
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = "D"
$ws.Range("L101").Value = 1.25
$ws.Range("M101").Value = 5.5
$ws.Range("N101").Value = 7.5
$ws.Range("O101").Value = 1.45
$ws.Range("P101").Value = 4.5
$ws.Range("Q101").Value = 5
$ws.Range("R101").Value = -1
$ws.Range("S101").Value = 1.775
$ws.Range("T101").Value = 2.025
$ws.Range("U101").Value = 2.5
$ws.Range("V101").Value = 1.875
$ws.Range("W101").Value = 1.925
$ws.Range("X101").Value = -1
$ws.Range("Y101").Value = 3.5
$ws.Range("Z101").Value = -1
$ws.Range("AA101").Value = -1
$ws.Range("AB101").Value = 1.025
$ws.Range("AC101").Value = -1
$ws.Range("AD101").Value = 0.925

# Row 102
$ws.Range("B102").Value = 6732837
$ws.Range("E102").Value = "Suduva Marijampole"
$ws.Range("F102").Value = "FK Riteriai"
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 3
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 1
$ws.Range("K102").Value = "A"
$ws.Range("L102").Value = 3.6
$ws.Range("M102").Value = 3.6
$ws.Range("N102").Value = 1.8
$ws.Range("O102").Value = 3
$ws.Range("P102").Value = 3.6
$ws.Range("Q102").Value = 2
$ws.Range("R102").Value = 0.25
$ws.Range("S102").Value = 2
$ws.Range("T102").Value = 1.8
$ws.Range("U102").Value = 2.5
$ws.Range("V102").Value = 1.975
$ws.Range("W102").Value = 1.825
$ws.Range("X102").Value = -1
$ws.Range("Y102").Value = -1
$ws.Range("Z102").Value = 1
$ws.Range("AA102").Value = -1
$ws.Range("AB102").Value = 0.8
$ws.Range("AC102").Value = 0.9750000000000001
$ws.Range("AD102").Value = -1

# Row 103
$ws.Range("B103").Value = 7465686
$ws.Range("E103").Value = "FK Kauno Zalgiris"
$ws.Range("F103").Value = "Hegelmann Litauen"
$ws.Range("G103").Value = 4
$ws.Range("H103").Value = 2
$ws.Range("I103").Value = 2
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = "H"
$ws.Range("L103").Value = 2.3
$ws.Range("M103").Value = 4
$ws.Range("N103").Value = 2.3
$ws.Range("O103").Value = 2.55
$ws.Range("P103").Value = 4
$ws.Range("Q103").Value = 2.2
$ws.Range("R103").Value = 0.25
$ws.Range("S103").Value = 1.8
$ws.Range("T103").Value = 2
$ws.Range("U103").Value = 2.75
$ws.Range("V103").Value = 1.85
$ws.Range("W103").Value = 1.95
$ws.Range("X103").Value = 1.55
$ws.Range("Y103").Value = -1
$ws.Range("Z103").Value = -1
$ws.Range("AA103").Value = 0.8
$ws.Range("AB103").Value = -1
$ws.Range("AC103").Value = 0.8500000000000001
$ws.Range("AD103").Value = -1

# Row 104
$ws.Range("B104").Value = 6732727
$ws.Range("E104").Value = "FK Zalgiris Vilnius"
$ws.Range("F104").Value = "FK Dainava Alytus"
$ws.Range("G104").Value = 1
$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 1
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = "H"
$ws.Range("L104").Value = 1.285
$ws.Range("M104").Value = 5.5
$ws.Range("N104").Value = 6.5
$ws.Range("O104").Value = 1.3
$ws.Range("P104").Value = 5.5
$ws.Range("Q104").Value = 6
$ws.Range("R104").Value = -1.5
$ws.Range("S104").Value = 1.9
$ws.Range("T104").Value = 1.9
$ws.Range("U104").Value = 2.75
$ws.Range("V104").Value = 1.8
$ws.Range("W104").Value = 2
$ws.Range("X104").Value = 0.3
$ws.Range("Y104").Value = -1
$ws.Range("Z104").Value = -1
$ws.Range("AA104").Value = -1
$ws.Range("AB104").Value = 0.8999999999999999
$ws.Range("AC104").Value = -1
$ws.Range("AD104").Value = 1

# Row 117
$ws.Range("B117").Value = 7862911
$ws.Range("E117").Value = "Hegelmann Litauen"
$ws.Range("F117").Value = "FK Siauliai"
$ws.Range("G117").Value = 2
$ws.Range("H117").Value = 2
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = "D"
$ws.Range("L117").Value = 2.15
$ws.Range("M117").Value = 3.1
$ws.Range("N117").Value = 3.1
$ws.Range("O117").Value = 2.45
$ws.Range("P117").Value = 2.9
$ws.Range("Q117").Value = 3
$ws.Range("R117").Value = 0
$ws.Range("S117").Value = 1.725
$ws.Range("T117").Value = 2.075
$ws.Range("U117").Value = 2.5
$ws.Range("V117").Value = 2.025
$ws.Range("W117").Value = 1.775
$ws.Range("X117").Value = -1
$ws.Range("Y117").Value = 1.9
$ws.Range("Z117").Value = -1
$ws.Range("AA117").Value = 0
$ws.Range("AB117").Value = 0
$ws.Range("AC117").Value = 1.025
$ws.Range("AD117").Value = -1

# Row 118
$ws.Range("B118").Value = 7862036
$ws.Range("E118").Value = "Banga Gargzdai"
$ws.Range("F118").Value = "FK Zalgiris Vilnius"
$ws.Range("G118").Value = 1
$ws.Range("H118").Value = 4
$ws.Range("I118").Value = 1
$ws.Range("J118").Value = 2
$ws.Range("K118").Value = "A"
$ws.Range("L118").Value = 8
$ws.Range("M118").Value = 4.5
$ws.Range("N118").Value = 1.3
$ws.Range("O118").Value = 6.5
$ws.Range("P118").Value = 4.5
$ws.Range("Q118").Value = 1.333
$ws.Range("R118").Value = 1.25
$ws.Range("S118").Value = 2
$ws.Range("T118").Value = 1.8
$ws.Range("U118").Value = 2.5
$ws.Range("V118").Value = 1.825
$ws.Range("W118").Value = 1.975
$ws.Range("X118").Value = -1
$ws.Range("Y118").Value = -1
$ws.Range("Z118").Value = 0.333
$ws.Range("AA118").Value = -1
$ws.Range("AB118").Value = 0.8
$ws.Range("AC118").Value = 1.025
$ws.Range("AD118").Value = -1
